$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Stage a copy of the original row 2 (values + formats) in a scratch row (row 20)
# so that the in-place rotation below doesn't clobber source data before it's read.
$ws.Range("A2:H2").Copy()
$ws.Range("A20:H20").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A2:H2").Copy()
$ws.Range("A20:H20").PasteSpecial(-4163)  # xlPasteValues

# The row-2 cells are rotated two columns to the left (C2->A2, D2->B2, ..., A2->G2, B2->H2).
# For each destination column (1=A .. 8=H), the matching staged source column is (destCol+2)
# wrapping around after column 8.
$srcCols = @(3, 4, 5, 6, 7, 8, 1, 2)

for ($i = 0; $i -lt 8; $i++) {
    $destCol = $i + 1
    $srcCol = $srcCols[$i]
    $srcCell = $ws.Cells.Item(20, $srcCol)
    $destCell = $ws.Cells.Item(2, $destCol)

    $srcCell.Copy()
    $destCell.PasteSpecial(-4122)  # xlPasteFormats
    $srcCell.Copy()
    $destCell.PasteSpecial(-4163)  # xlPasteValues
}

# Remove the scratch row used for staging
$ws.Range("A20:H20").Clear()

# Update the selection to match the new active cell / selection range
$ws.Range("G2:H2").Select()
